{"js": "// Implements the \"Aug 13 2023 feedback\" wording change in the\n// numbered instructions table of the Power of Attorney for Health\n// Care instructions document:\n//\n//   \"Sign and date the document. ... Your witness or witnesses has to\n//   see you sign it, and they must sign it right after you do.\"\n//\n// becomes\n//\n//   \"Sign and date the document in front of a witness. ... Your\n//   witness must see you sign it, and they must sign it right after\n//   you do.\"\n\nasync function replaceFirstMatch(searchText, replacementText, options) {\n  const searchOptions = Object.assign({ matchCase: true, matchWholeWord: false }, options || {});\n  const results = context.document.body.search(searchText, searchOptions);\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    return false;\n  }\n\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n  return true;\n}\n\n// 1) \"document.\" -> \"document in front of a witness.\"\nawait replaceFirstMatch(\"Sign and date the document.\", \"Sign and date the document in front of a witness.\");\n\n// 2) drop the \"or witnesses\" alternative - only \"witness\" remains\nawait replaceFirstMatch(\" or witnesses\", \"\");\n\n// 3) \"has\" -> \"must\" (whole word, so \"has\" inside other words isn't touched)\nawait replaceFirstMatch(\"has\", \"must\", { matchWholeWord: true });\n\n// 4) \"... must to see you sign it\" -> \"... must see you sign it\"\nawait replaceFirstMatch(\" to see you sign it\", \" see you sign it\");\n", "ps1": "# Implements the \"Aug 13 2023 feedback\" wording change in the\n# numbered instructions table of the Power of Attorney for Health Care\n# instructions document:\n#\n#   \"Sign and date the document. ... Your witness or witnesses has to\n#   see you sign it, and they must sign it right after you do.\"\n#\n# becomes\n#\n#   \"Sign and date the document in front of a witness. ... Your\n#   witness must see you sign it, and they must sign it right after\n#   you do.\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-InParagraph {\n    param(\n        [int]$ParaIndex,\n        [string]$SearchText,\n        [string]$ReplaceText,\n        [bool]$WholeWord = $false\n    )\n\n    $para = $d.Paragraphs.Item($ParaIndex)\n    $scope = $d.Range($para.Range.Start, $para.Range.End)\n\n    $find = $scope.Find\n    $find.ClearFormatting()\n    $find.Text = $SearchText\n    $find.MatchWholeWord = $WholeWord\n    $find.MatchCase = $true\n    $find.Wrap = 0          # wdFindStop - stay inside this paragraph only\n\n    $found = $find.Execute()\n    if ($found) {\n        $scope.Text = $ReplaceText\n    }\n    return $found\n}\n\n# The \"Sign and date the document...\" bullet is the 13th paragraph in\n# the document (the table cell holding it is its own paragraph).\n$signParaIndex = 13\n\n# 1) \"document.\" -> \"document in front of a witness.\"\nReplace-InParagraph $signParaIndex \"document.\" \"document in front of a witness.\" $false | Out-Null\n\n# 2) drop the \"or witnesses\" alternative - only \"witness\" remains\nReplace-InParagraph $signParaIndex \" or witnesses\" \"\" $false | Out-Null\n\n# 3) \"has\" -> \"must\"\nReplace-InParagraph $signParaIndex \"has\" \"must\" $true | Out-Null\n\n# 4) \"... must to see you sign it\" -> \"... must see you sign it\"\nReplace-InParagraph $signParaIndex \" to see you sign it\" \" see you sign it\" $false | Out-Null\n"}
